# Apply updated crypto price / % change figures per the latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.629.17"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "2.665.74"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'598.14"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("D6").Value = "'175.94"
$ws.Range("E6").Value = "  -1.78%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("D9").Value = "2.665.39"
$ws.Range("E9").Value = "  +1.22%  "
$ws.Range("E10").Value = "  -3.39%  "
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("D12").Value = "'0.357"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("D14").Value = "3.155.60"
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("D15").Value = "'0.0000186"
$ws.Range("E15").Value = "  -2.26%  "
$ws.Range("D16").Value = "72.456.11"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("E17").Value = "  -2.00%  "
$ws.Range("D18").Value = "2.658.81"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").Value = "'12.40"
$ws.Range("E19").Value = "  +5.12%  "
$ws.Range("D20").Value = "'8.20"
$ws.Range("E20").Value = "  +3.41%  "
$ws.Range("D21").Value = "'371.83"
$ws.Range("E21").Value = "  -3.60%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  +2.06%  "
$ws.Range("D24").Value = "'72.05"
$ws.Range("E24").Value = "  -2.84%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -1.68%  "
$ws.Range("D27").Value = "'9.86"
$ws.Range("E27").Value = "  -1.46%  "
$ws.Range("E28").Value = "  +2.16%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("D32").Value = "'497.49"
$ws.Range("E32").Value = "  -4.63%  "
$ws.Range("D33").Value = "'1.31"
$ws.Range("E33").Value = "  -2.11%  "
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("D36").Value = "'162.17"
$ws.Range("E36").Value = "  -1.56%  "
$ws.Range("D37").Value = "'19.54"
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("E38").Value = "  +1.09%  "
$ws.Range("D39").Value = "'18.95"
$ws.Range("E39").Value = "  -0.80%  "
$ws.Range("E40").Value = "  -1.86%  "
$ws.Range("E41").Value = "  -4.85%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").Value = "'5.02"
$ws.Range("E43").Value = "  -2.85%  "
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").Value = "'157.03"
$ws.Range("E46").Value = "  +4.02%  "
$ws.Range("D47").Value = "'39.29"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("D49").Value = "'0.560"
$ws.Range("E49").Value = "  +2.55%  "
$ws.Range("E50").Value = "  +1.73%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0258"
$ws.Range("E51").Value = "  -3.07%  "
